# Add a new "db_errs" report row to the report_locations sheet,
# producing a db error report for cases that look strange.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("report_locations")

$ws.Range("B7").Value = "~/covid_case_reports/db_errs"
$ws.Range("A7").Value = "db_errs"
$ws.Range("C7:I7").Value = "X"

$ws.Range("B7").Select()
